$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a brand-new worksheet named "2022" before the current first
#    sheet (so the tab order becomes 2022, 2021, 2018).
# ---------------------------------------------------------------------
$ws2022 = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$ws2022.Name = "2022"

# ---------------------------------------------------------------------
# 2. Populate the new "2022" sheet with the lab members / addresses for
#    that year (same two-column layout as the other year sheets).
# ---------------------------------------------------------------------
$data2022 = @(
    @("Name", "Address"),
    @("Ido Bar", "Givat Ada, Israel"),
    @("Sam (Prabhakaran) Thanjavur Sambasivam", "Chennai, Tamilnadu, India"),
    @("Chat Kanchana-Udomkan", "Kanchanaburi, Thailand"),
    @("Lab", "Nathan, QLD, Australia"),
    @("Rebecca Ford", "Hartley Wintney, Hampshire, UK"),
    @("Shane (Zi Wei) Zhou", "Nantong, Jiangsu, China"),
    @("Rachel Diva Soh", "Jakarta, Indonesia"),
    @("Melody Christie", "Prince Edward Island, Canada"),
    @("Jeremy Brownlie", "Canberra, Australia"),
    @("Gurpreet Singh Khalsa", "Jammu, India"),
    @("Matin Ghaheri", "Isfahan, Iran "),
    @("Henrietta Soi", "Nairobi, Kenya"),
    @("Fawad Ali", "Faisalabad-Punjab Pakistan"),
    @("Moutoshi Chakraborty", "Naogaon, Bangladesh"),
    @("Hayley Wilson", "Melbourne and Tamworth, Australia"),
    @("Joshua Lomax", "Cairns, Autralia"),
    @("Mahmuda Binte Monsur", "Tangail, Dhaka, Bangladesh")
)

for ($i = 0; $i -lt $data2022.Length; $i++) {
    $row = $data2022[$i]
    $ws2022.Cells.Item($i + 1, 1).Value = $row[0]
    $ws2022.Cells.Item($i + 1, 2).Value = $row[1]
}

# ---------------------------------------------------------------------
# 3. Turn the new range into a table, matching the naming scheme already
#    used by the other year sheets (Table1 / Table13 -> Table14).
# ---------------------------------------------------------------------
$lastRow = $data2022.Length
$tbl2022 = $ws2022.ListObjects.Add(1, $ws2022.Range("A1:B$lastRow"), $null, 1)
$tbl2022.Name = "Table14"

# ---------------------------------------------------------------------
# 4. Add a threaded comment on the "Fawad Ali" row, left by Fawad Ali.
# ---------------------------------------------------------------------
$excel.UserName = "Fawad Ali"
$ws2022.Range("A14").AddCommentThreaded("Fawad Ali") | Out-Null

# ---------------------------------------------------------------------
# 5. Keep column sizing tidy on the new sheet (best-fit, like the other
#    author-maintained year sheets).
# ---------------------------------------------------------------------
$ws2022.Columns.Item(1).AutoFit() | Out-Null
$ws2022.Columns.Item(2).AutoFit() | Out-Null

# ---------------------------------------------------------------------
# 6. Restore per-sheet selections so the active cell/tab matches what a
#    user would see after finishing data entry on "2022".
# ---------------------------------------------------------------------
$ws2021 = $wb.Worksheets.Item("2021")
$ws2021.Range("B11").Select() | Out-Null

$ws2018 = $wb.Worksheets.Item("2018")
$ws2018.Range("A3:B3").Select() | Out-Null

$ws2022.Range("B$lastRow").Select() | Out-Null
